$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Columns.Item(8).ColumnWidth = 20
Write-Host "done"
